$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.120.07'
$ws.Range("E2").Value = '  -0.80%  '

$ws.Range("D3").Value = '2.946.85'
$ws.Range("E3").Value = '  -1.38%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = '374.98'
$ws.Range("E5").Value = '  -1.63%  '

$ws.Range("D6").Value = '101.43'
$ws.Range("E6").Value = '  -2.93%  '

$ws.Range("D7").Value = '0.538'
$ws.Range("E7").Value = '  -1.88%  '

$ws.Range("E8").Value = '  +0.07%  '

$ws.Range("D9").Value = '0.588'
$ws.Range("E9").Value = '  -1.48%  '

$ws.Range("D10").Value = '36.36'
$ws.Range("E10").Value = '  -2.55%  '

$ws.Range("E11").Value = '  -0.71%  '

$ws.Range("D12").Value = '0.0852'
$ws.Range("E12").Value = '  +0.44%  '

$ws.Range("D13").Value = '3.409.32'
$ws.Range("E13").Value = '  -1.25%  '

$ws.Range("E14").Value = '  -2.00%  '

$ws.Range("E15").Value = '  +0.08%  '

$ws.Range("D16").Value = '2.983.57'
$ws.Range("E16").Value = '  -0.17%  '

$ws.Range("E17").Value = '  +2.51%  '

$ws.Range("D18").Value = '11.00'
$ws.Range("E18").Value = '  +47.46%  '

$ws.Range("D19").Value = '51.079.04'
$ws.Range("E19").Value = '  -0.78%  '

$ws.Range("E20").Value = '  -6.77%  '

$ws.Range("D21").Value = '12.47'
$ws.Range("E21").Value = '  -3.73%  '

$ws.Range("D22").Value = '0.0₃0954'
$ws.Range("E22").Value = '  -1.07%  '

$ws.Range("D23").Value = '265.34'

$ws.Range("D24").Value = '68.68'
$ws.Range("E24").Value = '  -1.12%  '

$ws.Range("D25").Value = '3.15'
$ws.Range("E25").Value = '  +7.32%  '

$ws.Range("D26").Value = '8.20'
$ws.Range("E26").Value = '  -0.53%  '

$ws.Range("D27").Value = '7.61'
$ws.Range("E27").Value = '  -1.59%  '

$ws.Range("E28").Value = '  -0.01%  '

$ws.Range("B29").Value = 'Kaspa'
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D29").Value = '0.164'
$ws.Range("E29").Value = '  -3.88%  '

$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").Value = '25.62'
$ws.Range("E30").Value = '  -1.47%  '

$ws.Range("E31").Value = '  -5.04%  '

$ws.Range("E32").Value = '  +1.17%  '

$ws.Range("D33").Value = '50.84'
$ws.Range("E33").Value = '  -0.30%  '

$ws.Range("D34").Value = '2.06'
$ws.Range("E34").Value = '  -0.77%  '

$ws.Range("D35").Value = '33.42'
$ws.Range("E35").Value = '  -4.14%  '

$ws.Range("D36").Value = '0.0443'
$ws.Range("E36").Value = '  -1.88%  '

$ws.Range("E37").Value = '  -0.16%  '

$ws.Range("D38").Value = '3.17'
$ws.Range("E38").Value = '  +3.87%  '

$ws.Range("E39").Value = '  -1.04%  '

$ws.Range("D40").Value = '16.39'
$ws.Range("E40").Value = '  -4.75%  '

$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").Value = '2.49'
$ws.Range("E41").Value = '  -3.88%  '

$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").Value = '1.79'
$ws.Range("E42").Value = '  -2.92%  '

$ws.Range("D43").Value = '120.80'
$ws.Range("E43").Value = '  -1.32%  '

$ws.Range("D44").Value = '21.29'
$ws.Range("E44").Value = '  -2.60%  '

$ws.Range("D45").Value = '3.37'
$ws.Range("E45").Value = '  +2.92%  '

$ws.Range("E46").Value = '  -0.54%  '

$ws.Range("E47").Value = '  -2.37%  '

$ws.Range("E48").Value = '  -2.18%  '

$ws.Range("D49").Value = '1.992.61'
$ws.Range("E49").Value = '  -1.91%  '

$ws.Range("D50").Value = '0.0325'
$ws.Range("E50").Value = '  -2.26%  '

$ws.Range("E51").Value = '  +2.51%  '
